$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D contains numeric-looking text (e.g. thousand-dot formatted prices).
# Force text format while assigning so Excel does not coerce these into numbers,
# then restore the original (default) style so the cells keep their original formatting.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = '63.695.42'
$ws.Range("E2").Value = '  -0.94%  '
$ws.Range("D3").Value = '3.417.96'
$ws.Range("E3").Value = '  -2.45%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '579.24'
$ws.Range("E5").Value = '  -1.85%  '
$ws.Range("D6").Value = '128.97'
$ws.Range("E6").Value = '  -4.03%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  -1.51%  '
$ws.Range("E9").Value = '  +3.52%  '
$ws.Range("E10").Value = '  -0.11%  '
$ws.Range("D11").Value = '0.381'
$ws.Range("E11").Value = '  -1.20%  '
$ws.Range("D12").Value = '4.001.09'
$ws.Range("E12").Value = '  -2.43%  '
$ws.Range("E13").Value = '  -0.41%  '
$ws.Range("E14").Value = '  -2.64%  '
$ws.Range("D15").Value = '3.420.64'
$ws.Range("E15").Value = '  -2.34%  '
$ws.Range("D16").Value = '63.685.47'
$ws.Range("E16").Value = '  -0.99%  '
$ws.Range("D17").Value = '25.29'
$ws.Range("E17").Value = '  -1.40%  '
$ws.Range("D18").Value = '9.81'
$ws.Range("E18").Value = '  -0.53%  '
$ws.Range("D19").Value = '5.63'
$ws.Range("E19").Value = '  -2.00%  '
$ws.Range("D20").Value = '13.30'
$ws.Range("E20").Value = '  -1.48%  '
$ws.Range("D21").Value = '382.55'
$ws.Range("E21").Value = '  -2.64%  '
$ws.Range("D22").Value = '0.562'
$ws.Range("E22").Value = '  -1.55%  '
$ws.Range("D23").Value = '3.556.65'
$ws.Range("E23").Value = '  -2.41%  '
$ws.Range("D24").Value = '74.02'
$ws.Range("E24").Value = '  -0.81%  '
$ws.Range("E25").Value = '  +0.10%  '
$ws.Range("E26").Value = '  -4.95%  '
$ws.Range("D27").Value = '0.999'
$ws.Range("E27").Value = '  -0.13%  '
$ws.Range("D28").Value = '2.19'
$ws.Range("E28").Value = '  -2.79%  '
$ws.Range("D29").Value = '6.99'
$ws.Range("E29").Value = '  -4.94%  '
$ws.Range("D30").Value = '7.86'
$ws.Range("E30").Value = '  -4.31%  '
$ws.Range("E31").Value = '  -0.62%  '
$ws.Range("E32").Value = '  -4.13%  '
$ws.Range("D33").Value = '3.448.26'
$ws.Range("E33").Value = '  -2.24%  '
$ws.Range("E34").Value = '  -0.10%  '
$ws.Range("D35").Value = '22.74'
$ws.Range("E35").Value = '  -3.06%  '
$ws.Range("D36").Value = '5.12'
$ws.Range("E36").Value = '  -0.22%  '
$ws.Range("D37").Value = '6.71'
$ws.Range("E37").Value = '  -2.47%  '
$ws.Range("D38").Value = '164.00'
$ws.Range("E38").Value = '  -2.11%  '
$ws.Range("E39").Value = '  -2.49%  '
$ws.Range("D40").Value = '0.0767'
$ws.Range("E40").Value = '  -1.65%  '
$ws.Range("D41").Value = '0.783'
$ws.Range("E41").Value = '  -3.47%  '
$ws.Range("E42").Value = '  -0.05%  '
$ws.Range("D43").Value = '41.39'
$ws.Range("E43").Value = '  -1.00%  '
$ws.Range("D44").Value = '4.29'
$ws.Range("E44").Value = '  -2.38%  '
$ws.Range("E45").Value = '  -3.65%  '
$ws.Range("D46").Value = '23.26'
$ws.Range("E46").Value = '  -7.32%  '
$ws.Range("D47").Value = '1.10'
$ws.Range("E47").Value = '  -5.92%  '
$ws.Range("D48").Value = '6.69'
$ws.Range("E48").Value = '  -0.99%  '
$ws.Range("D49").Value = '0.885'
$ws.Range("E49").Value = '  -0.63%  '
$ws.Range("D50").Value = '2.276.12'
$ws.Range("E50").Value = '  -2.86%  '
$ws.Range("E51").Value = '  -2.42%  '

# Restore original style/format (no explicit number format) on column D cells.
$dRange.Style = "Normal"
